$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh reshuffles which historical record (Fecha/Variedad/
# Volumen/Precios/Unidad/Origen/etc.) lands on each data row. Mercado ID,
# Mercado, Region, Codreg, Categoria ID, Categoria and Clasificacion are
# identical for every row already, so only columns D,H,I,J,K,L,M,N,O,P,Q
# need to move. new row -> old row it now carries the data of:
$map = @{2=23; 3=7; 4=3; 5=22; 6=8; 7=14; 8=12; 9=15; 10=16; 11=25; 12=5; 13=10; 14=4; 15=9; 16=20; 17=11; 18=26; 19=27; 20=2; 21=17; 22=18; 23=19; 24=28; 25=21; 26=29; 27=6; 28=24; 29=13}

$cols = @("D","H","I","J","K","L","M","N","O","P","Q")

# Snapshot every source value first so overwrites don't clobber data that
# another row still needs to read.
$snapshot = @{}
for ($r = 2; $r -le 29; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($newRow in $map.Keys) {
    $oldRow = $map[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $src[$c]
    }
}
